$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 48, shifting existing rows 48:57 down to 49:58
$ws.Rows.Item(48).Insert()

# Fill in the values for the new row 48
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value = "La Araucanía"
$ws.Cells.Item(48, 4).Value = 44785
$ws.Cells.Item(48, 5).Value = 9
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100108
$ws.Cells.Item(48, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(48, 9).Value = 100108003
$ws.Cells.Item(48, 10).Value = "Maracuyá"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 35
$ws.Cells.Item(48, 14).Value = 36000
$ws.Cells.Item(48, 15).Value = 36000
$ws.Cells.Item(48, 16).Value = 36000
$ws.Cells.Item(48, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(48, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 19).Value = 2000
$ws.Cells.Item(48, 20).Value = 18
